$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that are formatted like numbers (e.g. "1.002",
# "237.55"); assigning such text through COM lets Excel auto-coerce them to
# real numbers. Force the destination cell to Text first so the literal
# string is preserved, then drop the explicit format back to Normal so the
# cell's style stays the same as every untouched cell around it.
function Set-TextValue($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$updates = @(
    @{ Row = 2;  D = "30.395.13";     E = "  -0.91%  " },
    @{ Row = 3;  D = "1.894.70";      E = "  +0.42%  " },
    @{ Row = 4;  D = "1.002";         E = "  -0.02%  " },
    @{ Row = 5;  D = "237.55";        E = "  +0.48%  " },
    @{ Row = 6;  D = "1.001";         E = "  -0.18%  " },
    @{ Row = 7;  D = "0.4902";        E = "  +0.48%  " },
    @{ Row = 8;  D = "0.2920";        E = "  +0.39%  " },
    @{ Row = 9;  D = "0.06681";       E = "  +0.13%  " },
    @{ Row = 10; D = "1.886.67";      E = "  +0.12%  " },
    @{ Row = 11; D = "16.87";         E = "  +1.35%  " },
    @{ Row = 12; D = "0.07336";       E = "  +1.15%  " },
    @{ Row = 13; D = "5.144";         E = "  +2.62%  " },
    @{ Row = 14; D = "87.27";         E = "  -2.13%  " },
    @{ Row = 15; D = "0.6632";        E = "  +0.49%  " },
    @{ Row = 16; D = "30.394.22";     E = "  -0.69%  " },
    @{ Row = 17; D = "13.40";         E = "  +3.35%  " },
    @{ Row = 18; D = "0.000007804";   E = "  -0.65%  " },
    @{ Row = 19; D = "0.9998" },
    @{ Row = 20; D = "2.143.61";      E = "  +0.40%  " },
    @{ Row = 21; D = "5.329";         E = "  +12.58%  " },
    @{ Row = 22; D = "1.001";         E = "  -0.18%  " },
    @{ Row = 23; D = "191.06";        E = "  +0.30%  " },
    @{ Row = 24; D = "6.089";         E = "  -0.02%  " },
    @{ Row = 25; D = "9.447";         E = "  +1.61%  " },
    @{ Row = 26; D = "163.18";        E = "  +2.43%  " },
    @{ Row = 27; D = "18.18";         E = "  -0.52%  " },
    @{ Row = 28; D = "1.933";         E = "  +5.38%  " },
    @{ Row = 29; D = "1.473";         E = "  +4.66%  " },
    @{ Row = 30; D = "4.326";         E = "  +1.96%  " },
    @{ Row = 31; D = "0.09175";       E = "  +1.93%  " },
    @{ Row = 32; D = "4.041";         E = "  +2.84%  " },
    @{ Row = 33; D = "0.05171";       E = "  +0.23%  " },
    @{ Row = 34; D = "0.7384";        E = "  +1.79%  " },
    @{ Row = 35; D = "1.097";         E = "  +1.60%  " },
    @{ Row = 36; D = "2.716";         E = "  +0.69%  " },
    @{ Row = 37; D = "0.01806";       E = "  -0.23%  " },
    @{ Row = 38; D = "2.675";         E = "  +0.33%  " },
    @{ Row = 39; D = "0.9220";        E = "  +0.10%  " },
    @{ Row = 40; D = "2.032";         E = "  -0.76%  " },
    @{ Row = 41; D = "0.4364";        E = "  -0.27%  " },
    @{ Row = 42; D = "5.910";         E = "  +3.43%  " },
    @{ Row = 43; D = "106.12";        E = "  +1.37%  " },
    @{ Row = 44; D = "0.9943";        E = "  -0.46%  " },
    @{ Row = 45; D = "68.91";         E = "  +20.68%  " },
    @{ Row = 46; D = "0.1362";        E = "  +2.58%  " },
    @{ Row = 47; D = "7.551";         E = "  +3.01%  " },
    @{ Row = 48; D = "9.022";         E = "  +4.07%  " },
    @{ Row = 49; D = "34.87";         E = "  +4.96%  " },
    @{ Row = 50; D = "0.05828";       E = "  -0.03%  " },
    @{ Row = 51; D = "0.3914";        E = "  -3.21%  " }
)

foreach ($u in $updates) {
    Set-TextValue $ws.Cells.Item($u.Row, 4) $u.D
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
